# "So ajuste de conteudo de exemplo" - just tweaking the sample/example
# contact data and the sheet's display name.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from the default "Plan1" to "contatos"
$ws.Name = "contatos"

# Update Atila Xavier's example phone numbers and e-mail (row 2)
$ws.Range("C2").Value = "021981130000"
$ws.Range("D2").Value = "021999747200"
$ws.Range("G2").Value = "atila.xx@gmail.com"

# Leave the active selection on G3, matching where the user ended up
$ws.Range("G3").Select()
